$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values recalculated (strikeouts -> K stat), per row
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
